$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for week ending 2021-02-07 (36 rows): Datum, Receipt Number, Konto, Beskrivning, Debet, Kredit
$newRows = @(
    @(44229, 6021019, 3011, "Order  6021019 Swish +46707393913", $null, 423.21),
    @(44229, 6021019, 2611, "Order  6021019 Swish +46707393913", $null, 50.79),
    @(44229, 6021019, 1930, "Order  6021019 Swish +46707393913", 474, $null),
    @(44229, 5021732, 3011, "Order 5021732 Swish +46764282407", $null, 801.79),
    @(44229, 5021732, 2611, "Order 5021732 Swish +46764282407", $null, 96.20999999999999),
    @(44229, 5021732, 1930, "Order 5021732 Swish +46764282407", 898, $null),
    @(44229, 9021920, 3011, "Order 9021920 Swish +46704008971", $null, 655.36),
    @(44229, 9021920, 2611, "Order 9021920 Swish +46704008971", $null, 78.64),
    @(44229, 9021920, 1930, "Order 9021920 Swish +46704008971", 734, $null),
    @(44231, 6040901, 3011, "Order 6040901 Swish +46736813550", $null, 1062.5),
    @(44231, 6040901, 2611, "Order 6040901 Swish +46736813550", $null, 127.5),
    @(44231, 6040901, 1930, "Order 6040901 Swish +46736813550", 1190, $null),
    @(44231, $null, 5460, "IKEA BARKARBY K0135", 1196, $null),
    @(44231, $null, 2641, "IKEA BARKARBY K0135", 299, $null),
    @(44231, $null, 1930, "IKEA BARKARBY K0135", $null, 1495),
    @(44231, $null, 6400, "FACEBK PR3NBYWY62 K6885", 257, $null),
    @(44231, $null, $null, "FACEBK PR3NBYWY62 K6885", 0, $null),
    @(44231, $null, 1930, "FACEBK PR3NBYWY62 K6885", $null, 257),
    @(44232, 6051243, 3011, "Order 6051243 Swish +46705757460", $null, 1239.29),
    @(44232, 6051243, 2611, "Order 6051243 Swish +46705757460", $null, 148.71),
    @(44232, 6051243, 1930, "Order 6051243 Swish +46705757460", 1388, $null),
    @(44232, 5051600, 3011, "Order 5051600 Swish +46793490885", $null, 928.5700000000001),
    @(44232, 5051600, 2611, "Order 5051600 Swish +46793490885", $null, 111.43),
    @(44232, 5051600, 1930, "Order 5051600 Swish +46793490885", 1040, $null),
    @(44232, $null, 4010, "M&S RB BROMMA K0135", 527.78, $null),
    @(44232, $null, 2645, "M&S RB BROMMA K0135", 63.33, $null),
    @(44232, $null, 1930, "M&S RB BROMMA K0135", $null, 591.11),
    @(44234, 1070903, 3011, "Order 1070903 Swish +46709224929", $null, 1008.93),
    @(44234, 1070903, 2611, "Order 1070903 Swish +46709224929", $null, 121.07),
    @(44234, 1070903, 1930, "Order 1070903 Swish +46709224929", 1130, $null),
    @(44234, 1070927, 3011, "Order 1070927 Swish +46723656673", $null, 806.25),
    @(44234, 1070927, 2611, "Order 1070927 Swish +46723656673", $null, 96.75),
    @(44234, 1070927, 1930, "Order 1070927 Swish +46723656673", 903, $null),
    @(44234, $null, 5670, "ST1 V#LLINGBY K0135", 733.41, $null),
    @(44234, $null, 2641, "ST1 V#LLINGBY K0135", 183.35, $null),
    @(44234, $null, 1930, "ST1 V#LLINGBY K0135", $null, 916.76)
)

$startRow = 104
$endRow = $startRow + $newRows.Count - 1

# Apply the date number format used by column A (style index reused from existing rows)
$dateRange = "A" + $startRow + ":A" + $endRow
$ws.Range($dateRange).NumberFormat = $ws.Range("A103").NumberFormat

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
